$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: the "Implement a Service layer ..." bullet is currently split
# across several runs (artifacts of tracked edits / copy-paste). Collapse it
# down to a single run containing the full sentence, leaving the paragraph
# mark / list formatting untouched.
# ---------------------------------------------------------------------------
$targetParaText = "Implement a Service layer between constructor and repository. This layer can be used to do validation."

$svcPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $targetParaText) {
        $svcPara = $p
        break
    }
}

if ($svcPara -ne $null) {
    $r = $svcPara.Range
    # Range over just the paragraph's text (exclude the trailing paragraph mark).
    $body = $d.Range($r.Start, $r.End - 1)

    # The runtime short-circuits a Range.Text assignment when the new text is
    # identical to the existing text (which it is here, just split over many
    # runs) - so nothing would actually merge into one run. Force an actual
    # content change first, then set the final text, which collapses every
    # run under the range into a single new run.
    $placeholder = "TEMP_PLACEHOLDER_MERGE_RUNS"
    $body.Text = $placeholder
    $merged = $d.Range($r.Start, $r.Start + $placeholder.Length)
    $merged.Text = $targetParaText
}

# ---------------------------------------------------------------------------
# Change 2: the final (empty) bulleted paragraph gets new body text.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
if ($lastPara.Range.Text.TrimEnd([char]13, [char]7) -eq "") {
    $lastPara.Range.Text = "Implement password confirmation and password recovery."
}
